$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 393: new positive-case count
$ws.Range("C393").Value = 92

# Row 394: new positive-case count + new hospital-admission count
$ws.Range("C394").Value = 84

# Column L is formatted as Text ("@"); write numerically like the rest of
# the column (General, then restore) instead of letting it coerce to a string.
$ws.Range("L394").NumberFormat = "General"
$ws.Range("L394").Value = 1
$ws.Range("L394").NumberFormat = "@"

# Row 395: fill in the day's figures (previously blank placeholders)
$ws.Range("C395").Value = 16
$ws.Range("E395").Value = 6
$ws.Range("F395").Value = 5
$ws.Range("G395").Value = 32

$ws.Range("L395").NumberFormat = "General"
$ws.Range("L395").Value = 0
$ws.Range("L395").NumberFormat = "@"

$ws.Range("M395").NumberFormat = "General"
$ws.Range("M395").Value = 0
$ws.Range("M395").NumberFormat = "@"

# Move the active selection on the frozen pane, as in the authored file
$ws.Range("O28").Select()
